# Scheduled runner update: refresh market-board derived values in Sheets
# (currentAveragePrice / LevePrice* / LeveProfit* columns) across all job sheets.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1849.3429
$ws.Range("I40").Value = 1495.76
$ws.Range("J40").Value = 2733.3
$ws.Range("K40").Value = 1495.76
$ws.Range("L40").Value = 2733.3
$ws.Range("M40").Value = -1320.76
$ws.Range("N40").Value = -3083.3
$ws.Range("H68").Value = 58117.8
$ws.Range("J68").Value = 58117.8
$ws.Range("L68").Value = 58117.8
$ws.Range("N68").Value = -59615.8
$ws.Range("H71").Value = 58117.8
$ws.Range("J71").Value = 58117.8
$ws.Range("L71").Value = 174353.4
$ws.Range("N71").Value = -181841.4
$ws.Range("H123").Value = 69717.11
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 69717.11
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 69717.11
$ws.Range("M123").ClearContents()
$ws.Range("N123").Value = -79517.11
$ws.Range("H132").Value = 2197.1428
$ws.Range("I132").Value = 2131.0637
$ws.Range("J132").Value = 3750
$ws.Range("K132").Value = 6393.1911
$ws.Range("L132").Value = 11250
$ws.Range("M132").Value = -3863.1911
$ws.Range("N132").Value = -16310
$ws.Range("H138").Value = 3453.3257
$ws.Range("J138").Value = 4227.8276
$ws.Range("L138").Value = 12683.4828
$ws.Range("N138").Value = -22963.4828

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 2112.3333
$ws.Range("I97").Value = 644.4286
$ws.Range("J97").Value = 7250
$ws.Range("K97").Value = 644.4286
$ws.Range("L97").Value = 7250
$ws.Range("M97").Value = -148.4286
$ws.Range("N97").Value = -8242
$ws.Range("H122").Value = 2768.889
$ws.Range("I122").Value = 2653.3333
$ws.Range("K122").Value = 7959.999899999999
$ws.Range("M122").Value = -5509.999899999999
$ws.Range("H132").Value = 4860.163
$ws.Range("I132").Value = 5031.9
$ws.Range("K132").Value = 15095.7
$ws.Range("M132").Value = -12565.7

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 36368.832
$ws.Range("I75").Value = 13106.5
$ws.Range("J75").Value = 48000
$ws.Range("K75").Value = 13106.5
$ws.Range("L75").Value = 48000
$ws.Range("M75").Value = -12170.5
$ws.Range("N75").Value = -49872
$ws.Range("H78").Value = 36368.832
$ws.Range("I78").Value = 13106.5
$ws.Range("J78").Value = 48000
$ws.Range("K78").Value = 39319.5
$ws.Range("L78").Value = 144000
$ws.Range("M78").Value = -34639.5
$ws.Range("N78").Value = -153360
$ws.Range("H94").Value = 573.95654
$ws.Range("I94").Value = 521.7143
$ws.Range("J94").Value = 1122.5
$ws.Range("K94").Value = 521.7143
$ws.Range("L94").Value = 1122.5
$ws.Range("M94").Value = -70.71429999999998
$ws.Range("N94").Value = -2024.5
$ws.Range("H105").Value = 142860960
$ws.Range("I105").Value = 250002860
$ws.Range("J105").Value = 5096.3335
$ws.Range("K105").Value = 250002860
$ws.Range("L105").Value = 5096.3335
$ws.Range("M105").Value = -250001113
$ws.Range("N105").Value = -8590.333500000001
$ws.Range("H134").Value = 3228.0637
$ws.Range("I134").Value = 2780.75
$ws.Range("K134").Value = 8342.25
$ws.Range("M134").Value = -5807.25

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 870.63336
$ws.Range("I22").Value = 652.0833
$ws.Range("J22").Value = 1016.3333
$ws.Range("K22").Value = 652.0833
$ws.Range("L22").Value = 1016.3333
$ws.Range("M22").Value = -302.0833
$ws.Range("N22").Value = -1716.3333
$ws.Range("H31").Value = 3340.3333
$ws.Range("I31").Value = 1239.8889
$ws.Range("J31").Value = 4390.5557
$ws.Range("K31").Value = 1239.8889
$ws.Range("L31").Value = 4390.5557
$ws.Range("M31").Value = -944.8888999999999
$ws.Range("N31").Value = -4980.5557
$ws.Range("H34").Value = 3340.3333
$ws.Range("I34").Value = 1239.8889
$ws.Range("J34").Value = 4390.5557
$ws.Range("K34").Value = 1239.8889
$ws.Range("L34").Value = 4390.5557
$ws.Range("M34").Value = -1037.8889
$ws.Range("N34").Value = -4794.5557
$ws.Range("H87").Value = 56545.453
$ws.Range("J87").Value = 56545.453
$ws.Range("L87").Value = 56545.453
$ws.Range("N87").Value = -58917.453
$ws.Range("H90").Value = 56545.453
$ws.Range("J90").Value = 56545.453
$ws.Range("L90").Value = 169636.359
$ws.Range("N90").Value = -181492.359
$ws.Range("H92").Value = 30269.4
$ws.Range("J92").Value = 30269.4
$ws.Range("L92").Value = 30269.4
$ws.Range("N92").Value = -35261.4

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 37261580
$ws.Range("I4").Value = 252233.05
$ws.Range("J4").Value = 143002560
$ws.Range("K4").Value = 756699.1499999999
$ws.Range("L4").Value = 429007680
$ws.Range("M4").Value = -756587.1499999999
$ws.Range("N4").Value = -429007904
$ws.Range("H68").Value = 1722.8387
$ws.Range("J68").Value = 1733.7273
$ws.Range("L68").Value = 5201.1819
$ws.Range("N68").Value = -6823.1819
$ws.Range("H71").Value = 1722.8387
$ws.Range("J71").Value = 1733.7273
$ws.Range("L71").Value = 15603.5457
$ws.Range("N71").Value = -23715.5457
$ws.Range("H80").Value = 12307.471
$ws.Range("I80").Value = 18093
$ws.Range("J80").Value = 8257.6
$ws.Range("K80").Value = 54279
$ws.Range("L80").Value = 24772.8
$ws.Range("M80").Value = -53343
$ws.Range("N80").Value = -26644.8
$ws.Range("H83").Value = 12307.471
$ws.Range("I83").Value = 18093
$ws.Range("J83").Value = 8257.6
$ws.Range("K83").Value = 162837
$ws.Range("L83").Value = 74318.40000000001
$ws.Range("M83").Value = -158157
$ws.Range("N83").Value = -83678.40000000001
$ws.Range("H117").Value = 290377.88
$ws.Range("J117").Value = 464084.34
$ws.Range("L117").Value = 1392253.02
$ws.Range("N117").Value = -1399137.02
$ws.Range("H137").Value = 1963.8
$ws.Range("I137").Value = 1547.5
$ws.Range("K137").Value = 4642.5
$ws.Range("M137").Value = 457.5

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H131").Value = 59271.5
$ws.Range("J131").Value = 59271.5
$ws.Range("L131").Value = 59271.5
$ws.Range("N131").Value = -69351.5

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2248.125
$ws.Range("I22").Value = 2637.8
$ws.Range("J22").Value = 1598.6666
$ws.Range("K22").Value = 2637.8
$ws.Range("L22").Value = 1598.6666
$ws.Range("M22").Value = -2342.8
$ws.Range("N22").Value = -2188.6666
$ws.Range("H27").Value = 2248.125
$ws.Range("I27").Value = 2637.8
$ws.Range("J27").Value = 1598.6666
$ws.Range("K27").Value = 2637.8
$ws.Range("L27").Value = 1598.6666
$ws.Range("M27").Value = -2530.8
$ws.Range("N27").Value = -1812.6666
$ws.Range("H40").Value = 7087
$ws.Range("I40").Value = 9451.5
$ws.Range("K40").Value = 9451.5
$ws.Range("M40").Value = -9315.5
$ws.Range("H55").Value = 899.0833
$ws.Range("I55").Value = 278.8
$ws.Range("K55").Value = 278.8
$ws.Range("M55").Value = -105.8
$ws.Range("H93").Value = 1882.7222
$ws.Range("I93").Value = 1799.3529
$ws.Range("J93").Value = 3300
$ws.Range("K93").Value = 1799.3529
$ws.Range("L93").Value = 3300
$ws.Range("M93").Value = -551.3529000000001
$ws.Range("N93").Value = -5796
$ws.Range("H122").Value = 3487.1904
$ws.Range("I122").Value = 3401.8823
$ws.Range("K122").Value = 10205.6469
$ws.Range("M122").Value = -7755.6469

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 37333.332
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 37333.332
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 37333.332
$ws.Range("M48").ClearContents()
$ws.Range("N48").Value = -38471.332
